# Add three new Eredivisie 2023-2024 match rows (rows 51, 52, 53) to the
# sheet, mirroring the existing layout in columns A:V.
# Row 50 is the last pre-existing data row; its formatting (bold/border
# style on column A, date/time number format on column E) is cloned onto
# the new rows via Copy + PasteSpecial(formats) so the new cells end up
# sharing the exact same cell styles as the rest of the sheet.

$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$xlPasteFormats = -4122

function Set-MatchRow {
    param($r, $idx, $home, $homeGoals, $away, $awayGoals, $matchDate, $homeOpenOdds, $homeOpenTime, $homeCloseOdds, $homeCloseTime, $drawOpenOdds, $drawOpenTime, $drawCloseOdds, $drawCloseTime, $awayOpenOdds, $awayOpenTime, $awayCloseOdds, $awayCloseTime, $url)

    $ws.Cells.Item($r, 1).Value = $idx
    $ws.Cells.Item($r, 2).Value = "netherlands"
    $ws.Cells.Item($r, 3).Value = "eredivisie"
    $ws.Cells.Item($r, 4).Value = "2023-2024"
    $ws.Cells.Item($r, 5).Value = $matchDate
    $ws.Cells.Item($r, 6).Value = $home
    $ws.Cells.Item($r, 7).Value = $homeGoals
    $ws.Cells.Item($r, 8).Value = $away
    $ws.Cells.Item($r, 9).Value = $awayGoals
    $ws.Cells.Item($r, 10).Value = $homeOpenOdds
    $ws.Cells.Item($r, 11).Value = $homeOpenTime
    $ws.Cells.Item($r, 12).Value = $homeCloseOdds
    $ws.Cells.Item($r, 13).Value = $homeCloseTime
    $ws.Cells.Item($r, 14).Value = $drawOpenOdds
    $ws.Cells.Item($r, 15).Value = $drawOpenTime
    $ws.Cells.Item($r, 16).Value = $drawCloseOdds
    $ws.Cells.Item($r, 17).Value = $drawCloseTime
    $ws.Cells.Item($r, 18).Value = $awayOpenOdds
    $ws.Cells.Item($r, 19).Value = $awayOpenTime
    $ws.Cells.Item($r, 20).Value = $awayCloseOdds
    $ws.Cells.Item($r, 21).Value = $awayCloseTime
    $ws.Cells.Item($r, 22).Value = $url

    # Column A uses the bold / bordered / centered style applied to every
    # "Indice" cell in the sheet; column E uses the date/time number
    # format applied to every "data_partida" cell. Clone both from the
    # last existing data row (50) so the new cells share the same style.
    $ws.Cells.Item(50, 1).Copy()
    $ws.Cells.Item($r, 1).PasteSpecial($xlPasteFormats)
    $ws.Cells.Item(50, 5).Copy()
    $ws.Cells.Item($r, 5).PasteSpecial($xlPasteFormats)
}

Set-MatchRow 51 50 "Ajax" 0 "Feyenoord" 4 45196.58333333334 3.25 "24/09/2023 14:29" 3.25 "24/09/2023 14:29" 3.96 "24/09/2023 14:29" 3.96 "24/09/2023 14:29" 2.14 "24/09/2023 14:29" 2.14 "24/09/2023 14:29" "https://www.betexplorer.com/football/netherlands/eredivisie/ajax-feyenoord/WxSagv71/"

Set-MatchRow 52 51 "PSV" 3 "G.A. Eagles" 0 45196.78125 1.22 "23/09/2023 19:12" 1.15 "27/09/2023 18:37" 7.56 "23/09/2023 19:12" 9.53 "27/09/2023 18:44" 11.97 "23/09/2023 19:12" 16.38 "27/09/2023 18:44" "https://www.betexplorer.com/football/netherlands/eredivisie/psv-g-a-eagles/nic4G49j/"

Set-MatchRow 53 52 "Twente" 1 "Vitesse" 0 45196.83333333334 1.59 "24/09/2023 16:13" 1.36 "27/09/2023 19:56" 4.43 "24/09/2023 16:13" 5.35 "27/09/2023 19:56" 5.28 "24/09/2023 16:13" 9.16 "27/09/2023 19:56" "https://www.betexplorer.com/football/netherlands/eredivisie/twente-vitesse/MHXlG1EE/"
